# ---------------------------------------------------------------------------
# Adds a new "2022-Q1" worksheet (holding-detail data) right before the
# "总计" (total) summary sheet, and updates the "总计" sheet with a new
# first data row summarizing the 2022-Q1 quarter.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet by copying the most recent quarterly
#    sheet ("2021-Q4"). Copying preserves the header/row styles (s="2")
#    used by that sheet, and places the new sheet immediately after it.
# ---------------------------------------------------------------------------
$prev = $wb.Worksheets.Item("2021-Q4")
$prev.Copy([System.Reflection.Missing]::Value, $prev)
$new = $wb.Worksheets.Item("2021-Q4 (2)")
$new.Name = "2022-Q1"

# The source sheet only had 6 data rows (rows 2-7); we need 11 (rows 2-12).
# Copy the formatting of the last existing data row down onto the new rows
# so column A keeps the same style as the rest of the table.
$new.Range("A7:H7").Copy()
$new.Range("A8:H12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Force the fund-code / numeric-looking text columns to be stored as text
# (otherwise values such as "009394" or "25.59" get silently converted to
# numbers).
$new.Range("B2:B12").NumberFormat = "@"
$new.Range("D2:G12").NumberFormat = "@"

$new.Range("A2").Value = 0
$new.Range("B2").Value = "161810"
$new.Range("C2").Value = "银华内需精选混合(LOF)"
$new.Range("D2").Value = "25.59"
$new.Range("E2").Value = "94.71"
$new.Range("F2").Value = "7.96"
$new.Range("G2").Value = "2.0370"
$new.Range("H2").Value = 4

$new.Range("A3").Value = 1
$new.Range("B3").Value = "009394"
$new.Range("C3").Value = "银华同力精选混合"
$new.Range("D3").Value = "20.03"
$new.Range("E3").Value = "94.68"
$new.Range("F3").Value = "6.30"
$new.Range("G3").Value = "1.2619"
$new.Range("H3").Value = 5

$new.Range("A4").Value = 2
$new.Range("B4").Value = "001302"
$new.Range("C4").Value = "前海开源金银珠宝主题精选混合A"
$new.Range("D4").Value = "8.61"
$new.Range("E4").Value = "91.91"
$new.Range("F4").Value = "7.83"
$new.Range("G4").Value = "0.6742"
$new.Range("H4").Value = 8

$new.Range("A5").Value = 3
$new.Range("B5").Value = "003304"
$new.Range("C5").Value = "前海开源沪港深核心资源灵活配置混合A"
$new.Range("D5").Value = "5.91"
$new.Range("E5").Value = "93.10"
$new.Range("F5").Value = "8.08"
$new.Range("G5").Value = "0.4775"
$new.Range("H5").Value = 5

$new.Range("A6").Value = 4
$new.Range("B6").Value = "002207"
$new.Range("C6").Value = "前海开源金银珠宝主题精选混合C"
$new.Range("D6").Value = "3.45"
$new.Range("E6").Value = "91.91"
$new.Range("F6").Value = "7.83"
$new.Range("G6").Value = "0.2701"
$new.Range("H6").Value = 8

$new.Range("A7").Value = 5
$new.Range("B7").Value = "900009"
$new.Range("C7").Value = "中信证券成长动力混合A"
$new.Range("D7").Value = "6.14"
$new.Range("E7").Value = "89.40"
$new.Range("F7").Value = "3.12"
$new.Range("G7").Value = "0.1916"
$new.Range("H7").Value = 6

$new.Range("A8").Value = 6
$new.Range("B8").Value = "003305"
$new.Range("C8").Value = "前海开源沪港深核心资源灵活配置混合C"
$new.Range("D8").Value = "2.19"
$new.Range("E8").Value = "93.10"
$new.Range("F8").Value = "8.08"
$new.Range("G8").Value = "0.1770"
$new.Range("H8").Value = 5

$new.Range("A9").Value = 7
$new.Range("B9").Value = "000663"
$new.Range("C9").Value = "国投瑞银美丽中国灵活配置混合"
$new.Range("D9").Value = "3.55"
$new.Range("E9").Value = "92.92"
$new.Range("F9").Value = "3.99"
$new.Range("G9").Value = "0.1416"
$new.Range("H9").Value = 9

$new.Range("A10").Value = 8
$new.Range("B10").Value = "011997"
$new.Range("C10").Value = "景顺长城安盈回报一年持有期混合型证券投资基金A"
$new.Range("D10").Value = "5.74"
$new.Range("E10").Value = "29.73"
$new.Range("F10").Value = "1.54"
$new.Range("G10").Value = "0.0884"
$new.Range("H10").Value = 6

$new.Range("A11").Value = 9
$new.Range("B11").Value = "900059"
$new.Range("C11").Value = "中信证券成长动力混合C"
$new.Range("D11").Value = "0.65"
$new.Range("E11").Value = "89.40"
$new.Range("F11").Value = "3.12"
$new.Range("G11").Value = "0.0203"
$new.Range("H11").Value = 6

$new.Range("A12").Value = 10
$new.Range("B12").Value = "011998"
$new.Range("C12").Value = "景顺长城安盈回报一年持有期混合型证券投资基金C"
$new.Range("D12").Value = "0.23"
$new.Range("E12").Value = "29.73"
$new.Range("F12").Value = "1.54"
$new.Range("G12").Value = "0.0035"
$new.Range("H12").Value = 6

# ---------------------------------------------------------------------------
# 2. Update the "总计" (total) summary sheet: insert a new row right below
#    the header for the 2022-Q1 quarter, pushing the older quarters down.
# ---------------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")
$tot.Rows.Item(2).Insert(-4121) | Out-Null   # xlShiftDown
$tot.Range("B2:D2").ClearFormats()

# restore the column-A index style (it is not carried over automatically)
$tot.Range("A3").Copy()
$tot.Range("A2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 11
$tot.Range("D2").Value = 5.34

# Renumber the index column for the rows that shifted down.
$tot.Range("A3").Value = 1
$tot.Range("A4").Value = 2
$tot.Range("A5").Value = 3
$tot.Range("A6").Value = 4
$tot.Range("A7").Value = 5

Write-Output "2022-Q1 sheet added and 总计 sheet updated"
